# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
#
# Applies the hardpoint-data update + cosmetic tab-color/selection changes
# described by the commit to both worksheets:
#   - S2LAR_Sedan_HambaLG_f  (sheet1 / front)
#   - S2LAR_Sedan_HambaLG_r  (sheet2 / rear)

$wb = $excel.ActiveWorkbook

$wsFront = $wb.Worksheets.Item("S2LAR_Sedan_HambaLG_f")
$wsRear  = $wb.Worksheets.Item("S2LAR_Sedan_HambaLG_r")

# ---------------------------------------------------------------------------
# 1. Tab colors: theme 7 / tint 0.39997558519241921 (Accent4, lighter)
#                 -> theme 8 / tint -0.249977111117893 (Accent5, darker)
#    Pre-resolved to the equivalent RGB Excel would render for the new
#    theme+tint combination (Accent5 #5B9BD5 shaded by -0.249977111117893).
# ---------------------------------------------------------------------------
$newTabColor = RGB(46, 117, 182)   # 2E75B6
$wsFront.Tab.Color = $newTabColor
$wsRear.Tab.Color = $newTabColor

# ---------------------------------------------------------------------------
# 2. Updated hardpoints - front sheet (S2LAR_Sedan_HambaLG_f)
# ---------------------------------------------------------------------------
$wsFront.Range("F23").Value = 0.15379999999999999
$wsFront.Range("G23").Value = 0.65
$wsFront.Range("H23").Value = 0.24

$wsFront.Range("G24").Value = 0.91
$wsFront.Range("H24").Value = 0.23

$wsFront.Range("G26").Value = 0.62
$wsFront.Range("H26").Value = 0.65
$wsFront.Range("G26:H26").NumberFormat = "0.00"

$wsFront.Range("G27").Value = 0.85
$wsFront.Range("H27").Value = 0.19
$wsFront.Range("G27:H27").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 3. Updated hardpoints - rear sheet (S2LAR_Sedan_HambaLG_r)
# ---------------------------------------------------------------------------
$wsRear.Range("F23").Value = 0.13
$wsRear.Range("G23").Value = 0.65
$wsRear.Range("H23").Value = 0.24

$wsRear.Range("F24").Value = 0.13
$wsRear.Range("G24").Value = 0.91
$wsRear.Range("H24").Value = 0.23

$wsRear.Range("F26").Value = 0.0026557142857142869
$wsRear.Range("G26").Value = 0.62
$wsRear.Range("H26").Value = 0.65
$wsRear.Range("G26:H26").NumberFormat = "0.00"

$wsRear.Range("F27").Value = -0.055166428571428582
$wsRear.Range("G27").Value = 0.85
$wsRear.Range("H27").Value = 0.19
$wsRear.Range("G27:H27").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 4. Selection / active-cell bookmarks left behind when the author saved.
# ---------------------------------------------------------------------------
$wsFront.Range("C39").Select()
$wsRear.Range("H39").Select()

# Leave the front sheet active/selected, matching tabSelected="1" in sheet1.
$wsFront.Activate()
